# Waiter page data: Hours Worked for the waiter record moves from 0 to 13,
# and the active selection ends up on G3 (the cell right below the edited one).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 13

$ws.Range("G3").Select()
